$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A9").Value = "sulfur"
$ws.Range("A10").Value = "potassium"
$ws.Range("A11").Value = "sodium molybdate"

$ws.Range("A12").Select()
